$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header column: TIMEZONE (P1), and the row-2 import/validation sample data.
# Shared-string insertion order matters (matches how the values were actually
# typed in by the author), so cells are written in this specific sequence.

$ws.Range("A2").Value = 4223532
$ws.Range("B2").Value = "fafaf"
$ws.Range("D2").Value = 52534
$ws.Range("L2").Value = "gg"
$ws.Range("M2").Value = "Perempuan"
$ws.Range("N2").Value = "SD"
$ws.Range("O2").Value = "hahaha"
$ws.Range("P1").Value = "TIMEZONE"
$ws.Range("P2").Value = "WITA"
$ws.Range("F2").Value = "gg1"
$ws.Range("G2").Value = "gg2"
$ws.Range("H2").Value = "gg3"

# Selection: row 2 header-through-data range, active cell on the last header.
$ws.Range("A1:J2").Select()
$ws.Range("J2").Activate()

# Window size (bookView) — best effort; engine may not expose this.
$win = $excel.ActiveWindow
$win.Width = 20490
$win.Height = 9630
